# Add Loader for ClassRoom
# Inserts two new course rows (HU104A / HU104B - Functional English, taught
# by Sharmeen) right after the existing "MT101B" row (row 3), pushing the
# rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 4, shifting existing data (old rows 4-58)
# down to rows 6-60.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(5).Insert()

# Populate the new rows with the new course/teacher data.
$ws.Range("A4").Value = "HU104A"
$ws.Range("B4").Value = "Functional English"
$ws.Range("C4").Value = "CS1A"
$ws.Range("D4").Value = "Sharmeen"

$ws.Range("A5").Value = "HU104B"
$ws.Range("B5").Value = "Functional English"
$ws.Range("C5").Value = "CS1B"
$ws.Range("D5").Value = "Sharmeen"

# Match the author's final selection (active cell on A5).
[void]$ws.Range("A5").Select()
